$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 13:46"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 8094879
$ws.Range("C4").Value = 4626
$ws.Range("D4").Value = 5227718
$ws.Range("E4").Value = 2646222
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 66
$ws.Range("H4").Value = 220939

# Row 16: Iran -> Iran
$ws.Range("B16").Value = 513219
$ws.Range("C16").Value = 4830
$ws.Range("D16").Value = 414831
$ws.Range("E16").Value = 69039
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 279
$ws.Range("H16").Value = 29349

# Row 38: Catar -> Catar
$ws.Range("B38").Value = 128603
$ws.Range("C38").Value = 198
$ws.Range("D38").Value = 125584
$ws.Range("E38").Value = 2799
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 220

# Row 41: Nepal -> Nepal
$ws.Range("B41").Value = 117996
$ws.Range("C41").Value = 2638
$ws.Range("D41").Value = 80954
$ws.Range("E41").Value = 36367
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 12
$ws.Range("H41").Value = 675

# Row 70: Estado de Palestina -> Libia
$ws.Range("A70").Value = "Libia"
$ws.Range("B70").Value = 45821
$ws.Range("C70").Value = 836
$ws.Range("D70").Value = 25301
$ws.Range("E70").Value = 19851
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 13
$ws.Range("H70").Value = 669

# Row 71: Libia -> Estado de Palestina
$ws.Range("A71").Value = "Estado de Palestina"
$ws.Range("B71").Value = 45200
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 38841
$ws.Range("E71").Value = 5968
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 391

# Row 93: Malasia -> Malasia
$ws.Range("B93").Value = 17540
$ws.Range("C93").Value = 660
$ws.Range("D93").Value = 11605
$ws.Range("E93").Value = 5768
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 167

# Row 102: Finlandia -> Finlandia
$ws.Range("B102").Value = 12703
$ws.Range("C102").Value = 204
$ws.Range("D102").Value = 9100
$ws.Range("E102").Value = 3253
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 4
$ws.Range("H102").Value = 350

# Row 111: Luxemburgo -> Eslovenia
$ws.Range("A111").Value = "Eslovenia"
$ws.Range("B111").Value = 9938
$ws.Range("C111").Value = 707
$ws.Range("D111").Value = 5515
$ws.Range("E111").Value = 4248
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 175

# Row 112: Eslovenia -> Luxemburgo
$ws.Range("A112").Value = "Luxemburgo"
$ws.Range("B112").Value = 9840
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 8248
$ws.Range("E112").Value = 1459
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 133

# Row 139: Mayotte -> Malta
$ws.Range("A139").Value = "Malta"
$ws.Range("B139").Value = 4048
$ws.Range("C139").Value = 111
$ws.Range("D139").Value = 3064
$ws.Range("E139").Value = 940
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 44

# Row 140: Estonia -> Mayotte
$ws.Range("A140").Value = "Mayotte"
$ws.Range("B140").Value = 4030
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 2964
$ws.Range("E140").Value = 1023
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 43

# Row 141: Malta -> Estonia
$ws.Range("A141").Value = "Estonia"
$ws.Range("B141").Value = 3947
$ws.Range("C141").Value = 39
$ws.Range("D141").Value = 3060
$ws.Range("E141").Value = 819
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 68

# Row 143: Islandia -> Islandia
$ws.Range("B143").Value = 3757
$ws.Range("C143").Value = 89
$ws.Range("D143").Value = 2615
$ws.Range("E143").Value = 1132
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 10

# Row 168: Vietnam -> Vietnam
$ws.Range("B168").Value = 1122
$ws.Range("C168").Value = 9
$ws.Range("D168").Value = 1029
$ws.Range("E168").Value = 58
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 35

# Row 176: Taiwan -> Taiwan
$ws.Range("B176").Value = 530
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 491
$ws.Range("E176").Value = 32
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 7

# Row 178: Tanzania -> Gibraltar
$ws.Range("A178").Value = "Gibraltar"
$ws.Range("B178").Value = 516
$ws.Range("C178").Value = 17
$ws.Range("D178").Value = 435
$ws.Range("E178").Value = 81
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179: Gibraltar -> Tanzania
$ws.Range("A179").Value = "Tanzania"
$ws.Range("B179").Value = 509
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 183
$ws.Range("E179").Value = 305
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 21

# Row 181: Islas Feroe -> Islas Feroe
$ws.Range("B181").Value = 478
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 467
$ws.Range("E181").Value = 11
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

# Row 193: Liechtenstein -> Liechtenstein
$ws.Range("B193").Value = 174
$ws.Range("C193").Value = 11
$ws.Range("D193").Value = 131
$ws.Range("E193").Value = 42
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 1
